$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.406.64'
$ws.Range('E2').Value = '  -2.19%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.482.97'
$ws.Range('E3').Value = '  -2.20%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '610.12'
$ws.Range('E5').Value = '  +4.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '186.04'
$ws.Range('E6').Value = '  -0.29%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.628'
$ws.Range('E7').Value = '  -0.75%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  -1.98%  '
$ws.Range('E10').Value = '  -0.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.11'
$ws.Range('E11').Value = '  -2.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000308'
$ws.Range('E12').Value = '  -2.43%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.51'
$ws.Range('E13').Value = '  +0.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.030.33'
$ws.Range('E14').Value = '  -2.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '601.85'
$ws.Range('E15').Value = '  +5.62%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '69.389.87'
$ws.Range('E16').Value = '  -2.23%  '
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '12.62'
$ws.Range('E17').Value = '  +1.67%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.87'
$ws.Range('E18').Value = '  -1.75%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.481.00'
$ws.Range('E19').Value = '  -2.69%  '
$ws.Range('E20').Value = '  -0.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.986'
$ws.Range('E21').Value = '  -1.78%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.20'
$ws.Range('E22').Value = '  -2.61%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '105.99'
$ws.Range('E23').Value = '  +11.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.65'
$ws.Range('E24').Value = '  +2.18%  '
$ws.Range('E25').Value = '  +2.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.02'
$ws.Range('E26').Value = '  +2.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.95'
$ws.Range('E27').Value = '  -2.83%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.73'
$ws.Range('E28').Value = '  +6.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.57'
$ws.Range('E29').Value = '  +3.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.98'
$ws.Range('E30').Value = '  -3.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.43'
$ws.Range('E31').Value = '  +1.13%  '
$ws.Range('B32').Value = 'dogwifhat'
$ws.Range('C32').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.91'
$ws.Range('E32').Value = '  +15.62%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.116'
$ws.Range('E33').Value = '  -1.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.13'
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.18'
$ws.Range('E35').Value = '  -6.83%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '524.90'
$ws.Range('E37').Value = '  -4.29%  '
$ws.Range('E38').Value = '  -4.13%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.604.21'
$ws.Range('E39').Value = '  +1.15%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.62'
$ws.Range('E40').Value = '  +4.72%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.77'
$ws.Range('E41').Value = '  -2.98%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0₃0776'
$ws.Range('E42').Value = '  -3.41%  '
$ws.Range('E43').Value = '  +0.47%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0462'
$ws.Range('E44').Value = '  +0.88%  '
$ws.Range('E46').Value = '  +2.85%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.34'
$ws.Range('E47').Value = '  -4.23%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.81'
$ws.Range('E48').Value = '  -5.32%  '
$ws.Range('E49').Value = '  +0.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000243'
$ws.Range('E50').Value = '  -7.52%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '130.71'
$ws.Range('E51').Value = '  -3.32%  '
